# Updated CVDs for the month
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Milwaukee Pmc Hq Wisconsin (Professional Voluntary Turnover CVD update)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("E2").Value2 = 0.2326
$ws.Range("E3").Value2 = 0.2326
$ws.Range("E4").Value2 = 0.2326
$ws.Range("O4:W4").Value2 = 0
$ws.Range("O7").ClearContents()

# ---------------------------------------------------------------------
# Rosemont Illinois (Professional Voluntary Turnover CVD update)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rosemont Illinois")
$ws.Range("E2").Value2 = 0.0535
$ws.Range("E3").Value2 = 0.0535
$ws.Range("E4").Value2 = 0.0535
$ws.Range("O4:W4").Value2 = 0

# ---------------------------------------------------------------------
# Shanghai China (Professional Voluntary Turnover CVD update)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shanghai China")
$ws.Range("E2").Value2 = 0.2564
$ws.Range("E3").Value2 = 0.2564
$ws.Range("E4").Value2 = 0.2564
$ws.Range("O4:W4").Value2 = 0

# ---------------------------------------------------------------------
# St Barthelemy D'Anjou France (clear O4 Commit/Forecast cell)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("St Barthelemy D'Anjou France")
$ws.Range("O4").ClearContents()

# ---------------------------------------------------------------------
# Fort Wayne Indiana (clear O3 Commit/Forecast cell)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fort Wayne Indiana")
$ws.Range("O3").ClearContents()

# ---------------------------------------------------------------------
# Betzdorf Germany (Professional Voluntary Turnover CVD update)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Betzdorf Germany")
$ws.Range("E2").Value2 = 0.1316
$ws.Range("E3").Value2 = 0.1316
$ws.Range("E4").Value2 = 0.1316
$ws.Range("K4").Value2 = 0.125
$ws.Range("N4").Value2 = 0.137
$ws.Range("O4:W4").Value2 = 0

# ---------------------------------------------------------------------
# Hyderabad India (Professional Voluntary Turnover CVD update)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyderabad India")
$ws.Range("E2").Value2 = 0.2151
$ws.Range("E3").Value2 = 0.2151
$ws.Range("E4").Value2 = 0.2151
$ws.Range("O4").Value2 = 0.2222
